# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "RBF" source citation (rows 23-24) is being relocated further
# down the sheet to make room for a new "Definition of MSMEs" table. Remove
# the two rows outright (rather than just clearing values) so no stray
# formatting is left behind.
$ws.Rows("23:24").Delete()

# --- New "Definition of MSMEs" table (rows 19-23) ---------------------------
# Header row, styled like the other "title" header rows in the sheet (bold).
$ws.Range("B19").Value = "Number of employees"
$ws.Range("C19").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D19").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B19:D19").Font.Bold = $true

$ws.Range("A20").Value = "Micro"
$ws.Range("B20").Value = "<=5"

$ws.Range("A21").Value = "Small"
$ws.Range("B21").Value = "6-20"
$ws.Range("C21").Value = "USD 30,000 - 100,000"
$ws.Range("D21").Value = "USD 30,000 - 100,000"

$ws.Range("A22").Value = "Medium"
$ws.Range("B22").Value = "21-50"
$ws.Range("C22").Value = "USD 100,000 - 500,000"
$ws.Range("D22").Value = "USD 100,000 - 500,000"

$ws.Range("A23").Value = "Large"
$ws.Range("B23").Value = ">50"
$ws.Range("C23").Value = "USD > 500,000"
$ws.Range("D23").Value = "USD > 500,000"

# --- Relocated "RBF" source citation (rows 29-30) ---------------------------
$ws.Range("A29").Value = "RBF"
$ws.Range("A29").Font.Bold = $true

$ws.Range("A30").Value = "Reserve Bank of Fiji, ""ADDRESS BY MR INIA NAIYAGA, DEPUTY GOVERNOR OF THE RESERVE BANK OF FIJI AT THE FIJI INDIGENOUS BUSINESS COUNCIL 2014 SYMPOSIUM"", 2014, p.3. Available at http://www.rbf.gov.fj/docs2/Address%20by%20Mr%20Inia%20Naiyaga%20Deputy%20Governor%20of%20the%20Reserve%20Bank%20of%20Fiji%20at%20the%20Fiji%20Indigenous%20Business%20Council%202014%20Symposium.pdf"
$ws.Range("A30").Font.Italic = $true
